$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.840.15'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.735.44'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.00'
$ws.Range("E5").Value = '  +4.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06143'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '1.739.23'
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07166'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.92'
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6399'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.598'
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.04'
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D18").Value = '25.873.39'
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006753'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").Value = '1.962.90'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.261'
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.602'
$ws.Range("E23").Value = '  -2.05%  '
$ws.Range("E24").Value = '  +1.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.60'
$ws.Range("E25").Value = '  -1.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.514'
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.20'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.763'
$ws.Range("E28").Value = '  -2.27%  '
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.937'
$ws.Range("E30").Value = '  +5.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08237'
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.645'
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04628'
$ws.Range("E33").Value = '  +2.11%  '
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9854'
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.686'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01595'
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.918'
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.0000'
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.99'
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3830'
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7428'
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.992'
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.231'
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05237'
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.74'
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.548'
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3397'
$ws.Range("E51").Value = '  -1.01%  '
